$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as literal text, preserving the cell's
# existing style/number-format exactly (no new cellXfs entries) and without
# letting Excel auto-convert numeric-/date-looking strings into real numbers
# or dates. We compute the text through a formula in a scratch cell (so the
# copied value is guaranteed to be a text result, never re-parsed), then
# paste values-only onto the destination cell. Only used for cells that are
# not part of a merged range, since copy/paste-special would otherwise break
# the merge.
function Set-TextValue {
    param(
        [__ComObject]$Range,
        [string]$Text
    )
    $scratch = $ws.Range("ZZ1")
    $escaped = $Text -replace '"', '""'
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $Range.PasteSpecial(-4163) # xlPasteValues
    $scratch.Clear()
}

# ---- DATOS DEL USUARIO ----
# Plain text, safe to assign directly (keeps merges intact).
$ws.Range("C5").Value = "JORGE FRANKLYN COAQUIRA RAMOS"
$ws.Range("C6").Value = "jcoaquirar@pj.gob.pe"
# DNI looks numeric -> force text via scratch formula (G5 is not merged).
Set-TextValue $ws.Range("G5") "47087903"

# ---- FECHA / TIPO ACTA / FOLIO ----
# FECHA looks like a date -> force text (I3 is not merged).
Set-TextValue $ws.Range("I3") "02-01-2025"
$ws.Range("I4").Value = "ASIGNACIÓN"
$ws.Range("I7").Value = "7-2025"

# ---- BIEN (fila 14) ----
# CODIGO PATRIMONIAL looks numeric -> force text (B14 is not merged).
Set-TextValue $ws.Range("B14") "740800010013"
$ws.Range("D14").Value = "S/M"
$ws.Range("E14").Value = "S/M"
$ws.Range("F14").Value = "GRIS"
$ws.Range("G14").Value = "GFDFGH"

# ---- Firma (repite el nombre del usuario) ----
$ws.Range("F19").Value = "JORGE FRANKLYN COAQUIRA RAMOS"
